$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = -3.165797112063703
$ws.Range("F22").Value = -3.182967744055892
$ws.Range("F23").Value = -3.198819588169299
$ws.Range("F24").Value = -3.209031980779605
$ws.Range("F27").Value = -3.241847651373682
$ws.Range("F38").Value = -3.233693339851057
$ws.Range("F51").Value = -3.138536498339491
$ws.Range("F52").Value = -3.134143933150229
$ws.Range("F54").Value = -3.109606704160793
$ws.Range("F57").Value = -3.072530956742245
$ws.Range("F58").Value = -3.058796090317859
$ws.Range("F60").Value = -3.025560347044009
$ws.Range("F77").Value = -3.316652935775271
$ws.Range("F78").Value = -3.321028494846378
$ws.Range("F79").Value = -3.325836306140488
$ws.Range("F80").Value = -3.334729764760331
$ws.Range("F81").Value = -3.335830091148746
$ws.Range("F82").Value = -3.340086134142944
$ws.Range("F85").Value = -3.352382519657263
$ws.Range("F86").Value = -3.348572982404806
$ws.Range("F87").Value = -3.356371951108595
$ws.Range("F88").Value = -3.355858152877409
$ws.Range("F89").Value = -3.358744803095887
$ws.Range("F90").Value = -3.352326269161792
$ws.Range("F93").Value = -3.372443689809694
$ws.Range("F94").Value = -3.370006810279658
$ws.Range("F95").Value = -3.357193597181517
$ws.Range("F100").Value = -3.373489862685791
$ws.Range("F102").Value = -3.37454739320893
$ws.Range("F103").Value = -3.385644428882369
$ws.Range("F104").Value = -3.387908328663443
$ws.Range("F105").Value = -3.379126826505105
$ws.Range("F107").Value = -3.375501608988944
$ws.Range("F221").Value = -3.448550065400314
$ws.Range("F222").Value = -3.460445368484057
$ws.Range("F254").Value = -3.591092366674498
